# ----------------------------------------------------------------------
# Helper: resolve the 1-based Document.Paragraphs(..) index of whichever
# paragraph starts at a given character position. Paragraph.Index isn't
# reliable in this host, so we walk the collection and match on
# Range.Start instead.
# ----------------------------------------------------------------------
function Get-ParagraphIndexAt($doc, $startPos) {
    $n = $doc.Paragraphs.Count
    for ($i = 1; $i -le $n; $i++) {
        if ($doc.Paragraphs($i).Range.Start -eq $startPos) {
            return $i
        }
    }
    return -1
}

# Helper: find the paragraph whose text contains $needle and return its
# 1-based index in Document.Paragraphs.
function Find-ParagraphIndex($doc, $needle) {
    $scan = $doc.Content
    $scan.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $hitPara = $scan.Paragraphs(1)
    return Get-ParagraphIndexAt $doc $hitPara.Range.Start
}

$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Title text change. It occurs twice (under the "title" bookmark and
#    under the "firstheader" bookmark) — both literal, identical
#    Heading1 paragraphs. Re-find+replace each occurrence via the
#    paragraph's own Range (rather than Find's Replace) so Word keeps
#    xml:space="preserve" on the resulting run, matching a real edit.
# ------------------------------------------------------------------
$oldTitle = "Modelos de Factores Latentes en Econometría: Filtro de Kalman y Modelos de Espacio-Estado"
$newTitle = "Modelos de factores latentes dinámicos"

$n = $d.Paragraphs.Count
for ($i = 1; $i -le $n; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $oldTitle) {
        $p.Range.Text = $newTitle
    }
}

# ------------------------------------------------------------------
# 2. Abstract section: insert a new "AbstractFirstParagraph" paragraph
#    right before the "Palabras clave: ..." paragraph, carrying the new
#    placeholder abstract text, and demote the old paragraph's style
#    from FirstParagraph to BodyText.
# ------------------------------------------------------------------
$keywordsIdx = Find-ParagraphIndex $d "Palabras clave"

# InsertParagraphBefore() on paragraph N actually opens up a new blank
# paragraph immediately before paragraph N-1 in this host, so target
# paragraph N+1 (the one right after "Palabras clave") to land the new
# blank paragraph exactly where "Palabras clave" currently sits; the
# existing content then slides down by one.
$afterKeywords = $d.Paragraphs($keywordsIdx + 1)
$afterKeywords.Range.InsertParagraphBefore()

$newAbstractPara = $d.Paragraphs($keywordsIdx)
$newAbstractPara.Range.Text = "Este abstract será actualizado una vez que se complete el contenido final del artículo."
$newAbstractPara.Style = "AbstractFirstParagraph"

$keywordsPara = $d.Paragraphs($keywordsIdx + 1)
$keywordsPara.Style = "Body Text"

# ------------------------------------------------------------------
# 3. Remove the "Por Editar" list item from "Publicaciones Similares".
# ------------------------------------------------------------------
$porEditarIdx = Find-ParagraphIndex $d "Por Editar"
$d.Paragraphs($porEditarIdx).Range.Delete()

Write-Output "done"
